$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Cell value updates -----------------------------------------------
# Row 2
$ws.Range("B2").Value = "waseyrabby@btinternet.com"
$ws.Range("C2").Value = 718756

# Row 2 - Result column: clear full formatting then re-enter value (matches
# the observed loss of the cell style on N2 in the target workbook)
$ws.Range("N2").Clear()
$ws.Range("N2").Value = "Pass"

# Row 3
$ws.Range("B3").Value = "waseyrabby@btinternet.com"
$ws.Range("C3").Value = 718756
$ws.Range("N3").ClearContents()

# Row 4
$ws.Range("B4").Value = "waseyrabby@btinternet.com"

# Row 5
$ws.Range("A5").Value = "RegistrationTest"
$ws.Range("B5").Value = "waseyrabby@btinternet.com"

# Row 6
$ws.Range("B6").Value = "waseyrabby@btinternet.com"

# --- Column widths ------------------------------------------------------
$ws.Columns.Item(2).ColumnWidth = 25.833333333333332
$ws.Columns.Item(9).ColumnWidth = 35.833333333333336
$ws.Columns.Item(11).ColumnWidth = 16.166666666666668
$ws.Columns.Item(12).ColumnWidth = 17.0
$ws.Columns.Item(13).ColumnWidth = 28.333333333333332
$ws.Columns.Item(14).ColumnWidth = 25.5

# --- New hyperlinks -------------------------------------------------------
$ws.Hyperlinks.Add($ws.Range("C4"), "mailto:Test@123", [Type]::Missing, [Type]::Missing, "Test@123")
$ws.Hyperlinks.Add($ws.Range("C5"), "mailto:Test@123", [Type]::Missing, [Type]::Missing, "Test@123")
$ws.Hyperlinks.Add($ws.Range("B6"), "mailto:waseyrabby@btinternet.com")
$ws.Hyperlinks.Add($ws.Range("B2"), "mailto:waseyrabby@btinternet.com")
$ws.Hyperlinks.Add($ws.Range("B3"), "mailto:waseyrabby@btinternet.com")
$ws.Hyperlinks.Add($ws.Range("B4"), "mailto:waseyrabby@btinternet.com")
$ws.Hyperlinks.Add($ws.Range("B5"), "mailto:waseyrabby@btinternet.com")

# --- Page setup -----------------------------------------------------------
$ws.PageSetup.Orientation = 1

# --- View / selection -------------------------------------------------------
$excel.ActiveWindow.ScrollColumn = 2
$ws.Range("B5").Select()
